# Auto-generated Excel COM-interop script
# Applies the numeric cell updates described by the commit diff
# (scheduled market-data refresh across the Leve profit sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10237.333
$ws.Range("I64").Value = 7350
$ws.Range("K64").Value = 7350
$ws.Range("M64").Value = -7102
$ws.Range("H67").Value = 10237.333
$ws.Range("I67").Value = 7350
$ws.Range("K67").Value = 7350
$ws.Range("M67").Value = -6492
$ws.Range("H80").Value = 138889890
$ws.Range("I80").Value = 250000110
$ws.Range("K80").Value = 750000330
$ws.Range("M80").Value = -749999332
$ws.Range("H83").Value = 138889890
$ws.Range("I83").Value = 250000110
$ws.Range("K83").Value = 2250000990
$ws.Range("M83").Value = -2249995998
$ws.Range("H128").Value = 134274.28
$ws.Range("J128").Value = 134274.28
$ws.Range("L128").Value = 134274.28
$ws.Range("N128").Value = -144234.28
$ws.Range("H132").Value = 19347.418
$ws.Range("I132").Value = 23819.545
$ws.Range("K132").Value = 71458.63499999999
$ws.Range("M132").Value = -68928.63499999999
$ws.Range("H134").Value = 103747.25
$ws.Range("J134").Value = 103747.25
$ws.Range("L134").Value = 103747.25
$ws.Range("N134").Value = -113887.25
$ws.Range("H137").Value = 4001871.8
$ws.Range("I137").Value = 1773.3846
$ws.Range("K137").Value = 5320.1538
$ws.Range("M137").Value = -2770.1538
$ws.Range("H138").Value = 4106.1387
$ws.Range("I138").Value = 4655.25
$ws.Range("J138").Value = 3227.56
$ws.Range("K138").Value = 13965.75
$ws.Range("L138").Value = 9682.68
$ws.Range("M138").Value = -8825.75
$ws.Range("N138").Value = -19962.68
$ws.Range("H141").Value = 3685.8125
$ws.Range("I141").Value = 3133.8572
$ws.Range("K141").Value = 9401.571599999999
$ws.Range("M141").Value = -4221.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 8186
$ws.Range("J30").Value = 15997
$ws.Range("L30").Value = 15997
$ws.Range("N30").Value = -16297
$ws.Range("H74").Value = 772134.1
$ws.Range("I74").Value = 4064.8333
$ws.Range("K74").Value = 4064.8333
$ws.Range("M74").Value = -3190.8333
$ws.Range("H77").Value = 772134.1
$ws.Range("I77").Value = 4064.8333
$ws.Range("K77").Value = 20324.1665
$ws.Range("M77").Value = -15956.1665
$ws.Range("H122").Value = 1906.8
$ws.Range("I122").Value = 1805.0714
$ws.Range("J122").Value = 2313.7144
$ws.Range("K122").Value = 5415.2142
$ws.Range("L122").Value = 6941.1432
$ws.Range("M122").Value = -2965.2142
$ws.Range("N122").Value = -11841.1432
$ws.Range("H132").Value = 2187.077
$ws.Range("I132").Value = 492.44446
$ws.Range("K132").Value = 1477.33338
$ws.Range("M132").Value = 1052.66662
$ws.Range("H135").Value = 88456.39999999999
$ws.Range("J135").Value = 88456.39999999999
$ws.Range("L135").Value = 88456.39999999999
$ws.Range("N135").Value = -98596.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 12204.296
$ws.Range("I107").Value = 15512.158
$ws.Range("K107").Value = 15512.158
$ws.Range("M107").Value = -13592.158
$ws.Range("H134").Value = 21430512
$ws.Range("I134").Value = 1784.6765
$ws.Range("J134").Value = 112502610
$ws.Range("K134").Value = 5354.029500000001
$ws.Range("L134").Value = 337507830
$ws.Range("M134").Value = -2819.029500000001
$ws.Range("N134").Value = -337512900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1378.3966
$ws.Range("I58").Value = 1127.2766
$ws.Range("J58").Value = 2451.3635
$ws.Range("K58").Value = 1127.2766
$ws.Range("L58").Value = 2451.3635
$ws.Range("M58").Value = -924.2765999999999
$ws.Range("N58").Value = -2857.3635
$ws.Range("H134").Value = 1158.2142
$ws.Range("I134").Value = 1039.0435
$ws.Range("J134").Value = 1706.4
$ws.Range("K134").Value = 3117.1305
$ws.Range("L134").Value = 5119.200000000001
$ws.Range("M134").Value = -582.1305000000002
$ws.Range("N134").Value = -10189.2
$ws.Range("H136").Value = 1378.3966
$ws.Range("I136").Value = 1127.2766
$ws.Range("J136").Value = 2451.3635
$ws.Range("K136").Value = 3381.8298
$ws.Range("L136").Value = 7354.0905
$ws.Range("M136").Value = -831.8297999999995
$ws.Range("N136").Value = -12454.0905
$ws.Range("H140").Value = 161100.69
$ws.Range("J140").Value = 169109.08
$ws.Range("L140").Value = 169109.08
$ws.Range("N140").Value = -179469.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2855089.2
$ws.Range("I4").Value = 2578207
$ws.Range("J4").Value = 4101059.8
$ws.Range("K4").Value = 7734621
$ws.Range("L4").Value = 12303179.4
$ws.Range("M4").Value = -7734509
$ws.Range("N4").Value = -12303403.4
$ws.Range("H34").Value = 5701.4287
$ws.Range("J34").Value = 6583.3335
$ws.Range("L34").Value = 19750.0005
$ws.Range("N34").Value = -19918.0005
$ws.Range("H39").Value = 4247.1816
$ws.Range("I39").Value = 1057.3334
$ws.Range("J39").Value = 5443.375
$ws.Range("K39").Value = 3172.0002
$ws.Range("L39").Value = 16330.125
$ws.Range("M39").Value = -2878.0002
$ws.Range("N39").Value = -16918.125
$ws.Range("H58").Value = 15356.111
$ws.Range("I58").Value = 7401.6665
$ws.Range("J58").Value = 19333.334
$ws.Range("K58").Value = 22204.9995
$ws.Range("L58").Value = 58000.00199999999
$ws.Range("M58").Value = -22076.9995
$ws.Range("N58").Value = -58256.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1947
$ws.Range("I29").Value = 1899
$ws.Range("J29").Value = 1995
$ws.Range("K29").Value = 1899
$ws.Range("L29").Value = 1995
$ws.Range("M29").Value = -1609
$ws.Range("N29").Value = -2575
$ws.Range("H35").Value = 24999.5
$ws.Range("I35").Value = 24999.5
$ws.Range("K35").Value = 24999.5
$ws.Range("M35").Value = -24701.5
$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H70").Value = 4936.6665
$ws.Range("I70").Value = 4930
$ws.Range("K70").Value = 4930
$ws.Range("M70").Value = -4660
$ws.Range("H73").Value = 4936.6665
$ws.Range("I73").Value = 4930
$ws.Range("K73").Value = 4930
$ws.Range("M73").Value = -3994
$ws.Range("H80").Value = 52230144
$ws.Range("J80").Value = 45567044
$ws.Range("L80").Value = 45567044
$ws.Range("N80").Value = -45569040
$ws.Range("H83").Value = 52230144
$ws.Range("J83").Value = 45567044
$ws.Range("L83").Value = 227835220
$ws.Range("N83").Value = -227845204
$ws.Range("H102").Value = 27779238
$ws.Range("I102").Value = 33334686
$ws.Range("J102").Value = 1997.6666
$ws.Range("K102").Value = 33334686
$ws.Range("L102").Value = 1997.6666
$ws.Range("M102").Value = -33333064
$ws.Range("N102").Value = -5241.6666
$ws.Range("H132").Value = 646204.7
$ws.Range("I132").Value = 9107.643
$ws.Range("J132").Value = 989256.9399999999
$ws.Range("K132").Value = 27322.929
$ws.Range("L132").Value = 2967770.82
$ws.Range("M132").Value = -24792.929
$ws.Range("N132").Value = -2972830.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10905.083
$ws.Range("I7").Value = 6808.364
$ws.Range("K7").Value = 6808.364
$ws.Range("M7").Value = -6696.364
$ws.Range("H68").Value = 2474
$ws.Range("I68").Value = 2298.6667
$ws.Range("K68").Value = 2298.6667
$ws.Range("M68").Value = -1549.6667
$ws.Range("H71").Value = 2474
$ws.Range("I71").Value = 2298.6667
$ws.Range("K71").Value = 11493.3335
$ws.Range("M71").Value = -7749.333500000001
$ws.Range("H93").Value = 1053.625
$ws.Range("I93").Value = 801
$ws.Range("K93").Value = 801
$ws.Range("M93").Value = 447
$ws.Range("H126").Value = 10905.083
$ws.Range("I126").Value = 6808.364
$ws.Range("K126").Value = 20425.092
$ws.Range("M126").Value = -17955.092
$ws.Range("H132").Value = 3479.0688
$ws.Range("I132").Value = 3479.0688
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10437.2064
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -7907.206399999999
$ws.Range("H136").Value = 3924.5356
$ws.Range("I136").Value = 2994.8125
$ws.Range("J136").Value = 5164.1665
$ws.Range("K136").Value = 8984.4375
$ws.Range("L136").Value = 15492.4995
$ws.Range("M136").Value = -6434.4375
$ws.Range("N136").Value = -20592.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 14388.777
$ws.Range("J18").Value = 14333.167
$ws.Range("L18").Value = 14333.167
$ws.Range("N18").Value = -14679.167
$ws.Range("H80").Value = 34333
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 36499.5
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 36499.5
$ws.Range("M80").Value = -29002
$ws.Range("N80").Value = -38495.5
$ws.Range("H83").Value = 34333
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 36499.5
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 109498.5
$ws.Range("M83").Value = -85008
$ws.Range("N83").Value = -119482.5
$ws.Range("H98").Value = 57000
$ws.Range("J98").Value = 57000
$ws.Range("L98").Value = 57000
$ws.Range("N98").Value = -62990
$ws.Range("H132").Value = 1839.862
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

